# A new daily price record (Start Ruby / Primera, fecha 2021-09-28) was
# inserted into the price log right before the existing row for
# r=64 (fecha 2021-08-26). This pushes the existing rows 64-142 down to
# 65-143, and the sheet's used range grows from A1:T142 to A1:T143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 64, shifting rows 64:142 down to 65:143.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record.
$ws.Range("A64").Value = 10
$ws.Range("B64").Value = 'Vega Modelo de Temuco'
$ws.Range("C64").Value = 'La Araucanía'
$ws.Range("D64").Value = 44467
$ws.Range("E64").Value = 9
$ws.Range("F64").Value = 'Fruta'
$ws.Range("G64").Value = 100102
$ws.Range("H64").Value = 'Cítricos'
$ws.Range("I64").Value = 100102006
$ws.Range("J64").Value = 'Pomelo'
$ws.Range("K64").Value = 'Start Ruby'
$ws.Range("L64").Value = 'Primera'
$ws.Range("M64").Value = 80
$ws.Range("N64").Value = 10000
$ws.Range("O64").Value = 10000
$ws.Range("P64").Value = 10000
$ws.Range("Q64").Value = '$/bandeja 15 kilos granel'
$ws.Range("R64").Value = "Región de O'Higgins"
$ws.Range("S64").Value = 667
$ws.Range("T64").Value = 15
